$wb = $excel.ActiveWorkbook


# --- Sheet: ALC ---
$ws_ALC = $wb.Worksheets.Item("ALC")
$ws_ALC.Range("H12").Value = 1118.3334
$ws_ALC.Range("I12").Value = 3000
$ws_ALC.Range("J12").Value = 177.5
$ws_ALC.Range("K12").Value = 3000
$ws_ALC.Range("L12").Value = 177.5
$ws_ALC.Range("M12").Value = -2830
$ws_ALC.Range("N12").Value = -517.5
$ws_ALC.Range("H19").Value = 1714.75
$ws_ALC.Range("I19").Value = 3711.6667
$ws_ALC.Range("J19").Value = 516.6
$ws_ALC.Range("K19").Value = 3711.6667
$ws_ALC.Range("L19").Value = 516.6
$ws_ALC.Range("M19").Value = -3536.6667
$ws_ALC.Range("N19").Value = -866.6
$ws_ALC.Range("H43").Value = 439.8
$ws_ALC.Range("I43").Value = 0
$ws_ALC.Range("K43").Value = 0
$ws_ALC.Range("M43").ClearContents()
$ws_ALC.Range("H88").Value = 1271.6666
$ws_ALC.Range("I88").Value = 825.75
$ws_ALC.Range("J88").Value = 1628.4
$ws_ALC.Range("K88").Value = 825.75
$ws_ALC.Range("L88").Value = 1628.4
$ws_ALC.Range("M88").Value = -419.75
$ws_ALC.Range("N88").Value = -2440.4
$ws_ALC.Range("H91").Value = 1271.6666
$ws_ALC.Range("I91").Value = 825.75
$ws_ALC.Range("J91").Value = 1628.4
$ws_ALC.Range("K91").Value = 825.75
$ws_ALC.Range("L91").Value = 1628.4
$ws_ALC.Range("M91").Value = 578.25
$ws_ALC.Range("N91").Value = -4436.4
$ws_ALC.Range("H112").Value = 1134.4
$ws_ALC.Range("J112").Value = 1134.4
$ws_ALC.Range("L112").Value = 3403.2
$ws_ALC.Range("N112").Value = -5619.200000000001
$ws_ALC.Range("H113").Value = 40004264
$ws_ALC.Range("I113").Value = 100003780
$ws_ALC.Range("J113").Value = 4586.4
$ws_ALC.Range("K113").Value = 100003780
$ws_ALC.Range("L113").Value = 4586.4
$ws_ALC.Range("M113").Value = -100000526
$ws_ALC.Range("N113").Value = -11094.4
$ws_ALC.Range("H129").Value = 193874.14
$ws_ALC.Range("J129").Value = 197665.78
$ws_ALC.Range("L129").Value = 592997.34
$ws_ALC.Range("N129").Value = -602997.34
$ws_ALC.Range("H132").Value = 2564.275
$ws_ALC.Range("I132").Value = 2738.2222
$ws_ALC.Range("K132").Value = 8214.6666
$ws_ALC.Range("M132").Value = -5684.6666
$ws_ALC.Range("H138").Value = 1902.6666
$ws_ALC.Range("I138").Value = 877.2941
$ws_ALC.Range("J138").Value = 2203.2068
$ws_ALC.Range("K138").Value = 2631.8823
$ws_ALC.Range("L138").Value = 6609.6204
$ws_ALC.Range("M138").Value = 2508.1177
$ws_ALC.Range("N138").Value = -16889.6204

# --- Sheet: ARM ---
$ws_ARM = $wb.Worksheets.Item("ARM")
$ws_ARM.Range("H2").Value = 1488.4773
$ws_ARM.Range("I2").Value = 1428.4286
$ws_ARM.Range("J2").Value = 1722
$ws_ARM.Range("K2").Value = 1428.4286
$ws_ARM.Range("L2").Value = 1722
$ws_ARM.Range("M2").Value = -1315.4286
$ws_ARM.Range("N2").Value = -1948
$ws_ARM.Range("H32").Value = 25107.44
$ws_ARM.Range("I32").Value = 29883.146
$ws_ARM.Range("J32").Value = 3351.4443
$ws_ARM.Range("K32").Value = 29883.146
$ws_ARM.Range("L32").Value = 3351.4443
$ws_ARM.Range("M32").Value = -29596.146
$ws_ARM.Range("N32").Value = -3925.4443
$ws_ARM.Range("H59").Value = 22000
$ws_ARM.Range("J59").Value = 22000
$ws_ARM.Range("L59").Value = 22000
$ws_ARM.Range("N59").Value = -23608
$ws_ARM.Range("H97").Value = 881.6579
$ws_ARM.Range("I97").Value = 875.9
$ws_ARM.Range("K97").Value = 875.9
$ws_ARM.Range("M97").Value = -379.9
$ws_ARM.Range("H116").Value = 1488.4773
$ws_ARM.Range("I116").Value = 1428.4286
$ws_ARM.Range("J116").Value = 1722
$ws_ARM.Range("K116").Value = 1428.4286
$ws_ARM.Range("L116").Value = 1722
$ws_ARM.Range("M116").Value = 865.5714
$ws_ARM.Range("N116").Value = -6310
$ws_ARM.Range("H132").Value = 9361.937
$ws_ARM.Range("I132").Value = 1290.3962
$ws_ARM.Range("J132").Value = 52141.1
$ws_ARM.Range("K132").Value = 3871.188599999999
$ws_ARM.Range("L132").Value = 156423.3
$ws_ARM.Range("M132").Value = -1341.188599999999
$ws_ARM.Range("N132").Value = -161483.3

# --- Sheet: BSM ---
$ws_BSM = $wb.Worksheets.Item("BSM")
$ws_BSM.Range("H3").Value = 1488.4773
$ws_BSM.Range("I3").Value = 1428.4286
$ws_BSM.Range("J3").Value = 1722
$ws_BSM.Range("K3").Value = 1428.4286
$ws_BSM.Range("L3").Value = 1722
$ws_BSM.Range("M3").Value = -1314.4286
$ws_BSM.Range("N3").Value = -1950
$ws_BSM.Range("H20").Value = 2053
$ws_BSM.Range("I20").Value = 2723
$ws_BSM.Range("J20").Value = 936.3333
$ws_BSM.Range("K20").Value = 2723
$ws_BSM.Range("L20").Value = 936.3333
$ws_BSM.Range("M20").Value = -2476
$ws_BSM.Range("N20").Value = -1430.3333
$ws_BSM.Range("H92").Value = 27798.8
$ws_BSM.Range("J92").Value = 27798.8
$ws_BSM.Range("L92").Value = 27798.8
$ws_BSM.Range("N92").Value = -32790.8
$ws_BSM.Range("H128").Value = 2800
$ws_BSM.Range("I128").Value = 2800
$ws_BSM.Range("K128").Value = 8400
$ws_BSM.Range("M128").Value = -5910
$ws_BSM.Range("H134").Value = 36207.258
$ws_BSM.Range("I134").Value = 46351.125
$ws_BSM.Range("K134").Value = 139053.375
$ws_BSM.Range("M134").Value = -136518.375

# --- Sheet: CRP ---
$ws_CRP = $wb.Worksheets.Item("CRP")
$ws_CRP.Range("H31").Value = 10741.395
$ws_CRP.Range("I31").Value = 22634.934
$ws_CRP.Range("K31").Value = 22634.934
$ws_CRP.Range("M31").Value = -22339.934
$ws_CRP.Range("H34").Value = 10741.395
$ws_CRP.Range("I34").Value = 22634.934
$ws_CRP.Range("K34").Value = 22634.934
$ws_CRP.Range("M34").Value = -22432.934
$ws_CRP.Range("H58").Value = 17106.129
$ws_CRP.Range("I58").Value = 925.0476
$ws_CRP.Range("J58").Value = 51086.4
$ws_CRP.Range("K58").Value = 925.0476
$ws_CRP.Range("L58").Value = 51086.4
$ws_CRP.Range("M58").Value = -722.0476
$ws_CRP.Range("N58").Value = -51492.4
$ws_CRP.Range("H132").Value = 12541.872
$ws_CRP.Range("I132").Value = 14482.5
$ws_CRP.Range("J132").Value = 4348.1113
$ws_CRP.Range("K132").Value = 43447.5
$ws_CRP.Range("L132").Value = 13044.3339
$ws_CRP.Range("M132").Value = -40917.5
$ws_CRP.Range("N132").Value = -18104.3339
$ws_CRP.Range("H136").Value = 17106.129
$ws_CRP.Range("I136").Value = 925.0476
$ws_CRP.Range("J136").Value = 51086.4
$ws_CRP.Range("K136").Value = 2775.1428
$ws_CRP.Range("L136").Value = 153259.2
$ws_CRP.Range("M136").Value = -225.1428000000001
$ws_CRP.Range("N136").Value = -158359.2

# --- Sheet: CUL ---
$ws_CUL = $wb.Worksheets.Item("CUL")
$ws_CUL.Range("H5").Value = 973.6429000000001
$ws_CUL.Range("I5").Value = 439.36365
$ws_CUL.Range("J5").Value = 2932.6667
$ws_CUL.Range("K5").Value = 1318.09095
$ws_CUL.Range("L5").Value = 8798.000100000001
$ws_CUL.Range("M5").Value = -1206.09095
$ws_CUL.Range("N5").Value = -9022.000100000001
$ws_CUL.Range("H37").Value = 35755284
$ws_CUL.Range("J37").Value = 35755284
$ws_CUL.Range("L37").Value = 107265852
$ws_CUL.Range("N37").Value = -107266076
$ws_CUL.Range("H119").Value = 3666.6667
$ws_CUL.Range("I119").Value = 2400
$ws_CUL.Range("K119").Value = 7200
$ws_CUL.Range("M119").Value = -2362
$ws_CUL.Range("H131").Value = 796.5
$ws_CUL.Range("I131").Value = 250
$ws_CUL.Range("J131").Value = 813.4020400000001
$ws_CUL.Range("K131").Value = 750
$ws_CUL.Range("L131").Value = 2440.20612
$ws_CUL.Range("M131").Value = 4290
$ws_CUL.Range("N131").Value = -12520.20612
$ws_CUL.Range("H135").Value = 973.6429000000001
$ws_CUL.Range("I135").Value = 439.36365
$ws_CUL.Range("J135").Value = 2932.6667
$ws_CUL.Range("K135").Value = 3954.27285
$ws_CUL.Range("L135").Value = 26394.0003
$ws_CUL.Range("M135").Value = -1419.27285
$ws_CUL.Range("N135").Value = -31464.0003

# --- Sheet: GSM ---
$ws_GSM = $wb.Worksheets.Item("GSM")
$ws_GSM.Range("H102").Value = 26317046
$ws_GSM.Range("I102").Value = 31251220
$ws_GSM.Range("J102").Value = 1459.6666
$ws_GSM.Range("K102").Value = 31251220
$ws_GSM.Range("L102").Value = 1459.6666
$ws_GSM.Range("M102").Value = -31249598
$ws_GSM.Range("N102").Value = -4703.6666
$ws_GSM.Range("H126").Value = 3689.4443
$ws_GSM.Range("I126").Value = 3041.1304
$ws_GSM.Range("J126").Value = 4836.4614
$ws_GSM.Range("K126").Value = 9123.3912
$ws_GSM.Range("L126").Value = 14509.3842
$ws_GSM.Range("M126").Value = -6653.3912
$ws_GSM.Range("N126").Value = -19449.3842

# --- Sheet: LTW ---
$ws_LTW = $wb.Worksheets.Item("LTW")
$ws_LTW.Range("H22").Value = 10000
$ws_LTW.Range("I22").Value = 10000.5
$ws_LTW.Range("J22").Value = 9999.5
$ws_LTW.Range("K22").Value = 10000.5
$ws_LTW.Range("L22").Value = 9999.5
$ws_LTW.Range("M22").Value = -9705.5
$ws_LTW.Range("N22").Value = -10589.5
$ws_LTW.Range("H27").Value = 10000
$ws_LTW.Range("I27").Value = 10000.5
$ws_LTW.Range("J27").Value = 9999.5
$ws_LTW.Range("K27").Value = 10000.5
$ws_LTW.Range("L27").Value = 9999.5
$ws_LTW.Range("M27").Value = -9893.5
$ws_LTW.Range("N27").Value = -10213.5
$ws_LTW.Range("H55").Value = 71.22727
$ws_LTW.Range("J55").Value = 71.69231000000001
$ws_LTW.Range("L55").Value = 71.69231000000001
$ws_LTW.Range("N55").Value = -417.69231
$ws_LTW.Range("H61").Value = 4919.9165
$ws_LTW.Range("I61").Value = 3157.6
$ws_LTW.Range("J61").Value = 7857.1113
$ws_LTW.Range("K61").Value = 3157.6
$ws_LTW.Range("L61").Value = 7857.1113
$ws_LTW.Range("M61").Value = -2955.6
$ws_LTW.Range("N61").Value = -8261.1113
$ws_LTW.Range("H93").Value = 859.7059
$ws_LTW.Range("I93").Value = 859.7059
$ws_LTW.Range("K93").Value = 859.7059
$ws_LTW.Range("M93").Value = 388.2941
$ws_LTW.Range("H113").Value = 4919.9165
$ws_LTW.Range("I113").Value = 3157.6
$ws_LTW.Range("J113").Value = 7857.1113
$ws_LTW.Range("K113").Value = 3157.6
$ws_LTW.Range("L113").Value = 7857.1113
$ws_LTW.Range("M113").Value = -987.5999999999999
$ws_LTW.Range("N113").Value = -12197.1113
$ws_LTW.Range("H136").Value = 36709.285
$ws_LTW.Range("I136").Value = 50847
$ws_LTW.Range("K136").Value = 152541
$ws_LTW.Range("M136").Value = -149991

# --- Sheet: WVR ---
$ws_WVR = $wb.Worksheets.Item("WVR")
$ws_WVR.Range("H113").Value = 1229874.6
$ws_WVR.Range("I113").Value = 1906.5333
$ws_WVR.Range("J113").Value = 3861235
$ws_WVR.Range("K113").Value = 5719.5999
$ws_WVR.Range("L113").Value = 11583705
$ws_WVR.Range("M113").Value = -3549.5999
$ws_WVR.Range("N113").Value = -11588045
$ws_WVR.Range("H132").Value = 1336.9333
$ws_WVR.Range("I132").Value = 729.6667
$ws_WVR.Range("J132").Value = 3766
$ws_WVR.Range("K132").Value = 2189.0001
$ws_WVR.Range("L132").Value = 11298
$ws_WVR.Range("M132").Value = 340.9998999999998
$ws_WVR.Range("N132").Value = -16358
